$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Og endnu mere"
$ws.Range("C2").Value = 789

$ws.Range("C3").Select()
